# break out stock.yaml completed
#
# 1) Append rows 326-336 to the "day" sheet (new scrape batch, 12/08/2024).
# 2) Fix the "week" sheet so D132:D144 (bsecode) are stored as numbers
#    instead of text, matching the numeric-typed bsecode used elsewhere.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "day" sheet - append new rows
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$newDayRows = @(
    @(1,  "OFSS",       "Oracle Financial Services Software Limited",          532466, 1.61,                10734.55, 319423,   "day", "12/08/2024 11:35:41"),
    @(2,  "TCS",        "Tata Consultancy Services Limited",                   532540, -0.78,               4195.65,  936912,   "day", "12/08/2024 11:35:41"),
    @(3,  "METROPOLIS", "Metropolis Healthcare Ltd",                           542650, -0.9399999999999999, 2038.65,  236763,   "day", "12/08/2024 11:35:41"),
    @(4,  "INFY",       "Infosys Limited",                                     500209, 1.51,                1797.4,   4315329,  "day", "12/08/2024 11:35:41"),
    @(5,  "PVRINOX",    "PVR Inox Ltd",                                        532689, -1.74,               1475.1,   270377,   "day", "12/08/2024 11:35:41"),
    @(6,  "BHARTIARTL", "Bharti Airtel Limited",                               532454, -0.38,               1458.6,   4028226,  "day", "12/08/2024 11:35:41"),
    @(7,  "BATAINDIA",  "Bata India Limited",                                  500043, -1.68,               1418.8,   266322,   "day", "12/08/2024 11:35:41"),
    @(8,  "CHOLAFIN",   "Cholamandalam Investment And Finance Company Limited",511243, 0.01,                1348.75,  601506,   "day", "12/08/2024 11:35:41"),
    @(9,  "INDHOTEL",   "The Indian Hotels Company Limited",                   500850, -0.21,               617.15,   2938413,  "day", "12/08/2024 11:35:41"),
    @(10, "INDIACEM",   "The India Cements Limited",                           530005, 0.07000000000000001, 367.25,   1040771,  "day", "12/08/2024 11:35:41"),
    @(11, "PNB",        "Punjab National Bank",                                532461, -0.58,               114.6,    15811418, "day", "12/08/2024 11:35:41")
)

$startRow = 326
$r = $startRow
foreach ($row in $newDayRows) {
    $dayWs.Cells.Item($r, 1).Value = $row[0]
    $dayWs.Cells.Item($r, 2).Value = $row[1]
    $dayWs.Cells.Item($r, 3).Value = $row[2]
    $dayWs.Cells.Item($r, 4).Value = $row[3]
    $dayWs.Cells.Item($r, 5).Value = $row[4]
    $dayWs.Cells.Item($r, 6).Value = $row[5]
    $dayWs.Cells.Item($r, 7).Value = $row[6]
    $dayWs.Cells.Item($r, 8).Value = $row[7]
    $dayWs.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# ---------------------------------------------------------------------
# 2) "week" sheet - bsecode column (D) for rows 132-144 should be numeric
# ---------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

$bsecodes = @{
    132 = 540699
    133 = 532538
    134 = 500550
    135 = 500410
    136 = 532830
    137 = 500180
    138 = 500228
    139 = 511196
    140 = 500425
    141 = 512070
    142 = 500295
    143 = 517334
    144 = 533519
}

foreach ($rowNum in $bsecodes.Keys) {
    $weekWs.Cells.Item($rowNum, 4).Value = $bsecodes[$rowNum]
}

Write-Output "edit complete"
